$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.718.98"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "1.916.08"
$ws.Range("E3").Value = "  +0.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.96"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4930"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2979"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06792"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "1.898.89"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.24"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07365"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.183"
$ws.Range("E13").Value = "  +2.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.86"
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6751"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("D16").Value = "30.684.41"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007958"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.56"
$ws.Range("E18").Value = "  +3.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "2.187.23"
$ws.Range("E20").Value = "  +2.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.456"
$ws.Range("E21").Value = "  +13.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "201.91"
$ws.Range("E23").Value = "  +2.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.297"
$ws.Range("E24").Value = "  +2.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.682"
$ws.Range("E25").Value = "  +2.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.23"
$ws.Range("E26").Value = "  +5.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.74"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.962"
$ws.Range("E28").Value = "  +2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  +6.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.372"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09165"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05312"
$ws.Range("E33").Value = "  +1.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7424"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.712"
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01835"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("E38").Value = "  +1.11%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.103"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9226"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.54"
$ws.Range("E41").Value = "  +30.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4473"
$ws.Range("E42").Value = "  +1.06%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.951"
$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.22"
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1394"
$ws.Range("E46").Value = "  +2.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.667"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.178"
$ws.Range("E48").Value = "  +5.05%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.55"
$ws.Range("E49").Value = "  +5.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05874"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4039"
$ws.Range("E51").Value = "  +2.36%  "
